# Accommodate staging person unique ID; tweaks to dimensional summary tables
#
# This script:
#  1. Adds a new "ChargeTypeCitation" column (C) to the ChargeType sheet,
#     populating it with Citation 1-7 for the data rows and mirroring the
#     None/Unknown values used in column B for the trailer rows.
#  2. Makes the ChargeType sheet the active/selected sheet (it was previously
#     the MedicationType sheet that was active/selected), updating the
#     selection on both sheets accordingly.

$wb = $excel.ActiveWorkbook

$chargeType = $wb.Worksheets.Item("ChargeType")
$medicationType = $wb.Worksheets.Item("MedicationType")

# Header for the new column
$chargeType.Range("C1").Value = "ChargeTypeCitation"

# Citation 1 .. Citation 7 for rows 2-8
for ($i = 1; $i -le 7; $i++) {
    $row = $i + 1
    $chargeType.Cells.Item($row, 3).Value = "Citation $i"
}

# Trailer rows mirror column B ("None" / "Unknown") - use Value2 so we copy
# the plain text rather than a COM Variant wrapper.
$chargeType.Range("C9").Value = $chargeType.Range("B9").Value2
$chargeType.Range("C10").Value = $chargeType.Range("B10").Value2

# New column width (matches the other data columns' styling)
$chargeType.Columns.Item(3).ColumnWidth = 16.7

# Update selection on the sheet that used to be active (without activating it)
$medicationType.Range("E14").Select()

# Update selection on ChargeType and make it the active / selected tab
# (selecting a range on a sheet also activates that sheet, so do this last)
$chargeType.Range("C11").Select()
